# Fix reference error from ISO27002 12.7.1: the cross-reference in column E
# pointed to the nonexistent/duplicate "6.9.7.2" instead of "6.9.7.1".
# Also rename the "27001+27002" sheet to "ISO 27001+27002" for clarity.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("27001+27002")

# Rename the sheet.
$ws.Name = "ISO 27001+27002"

# Correct the mis-typed cross reference for ISO27002.12.7.1 (row 137, col E):
# it was "6.9.7.2" and should be "6.9.7.1".
$ws.Range("E137").Value = "6.9.7.1"

# Restore the sheet as the active/selected one, with the cursor parked on
# E138 (matches the author's on-screen selection after the fix).
$ws.Activate()
$ws.Range("E138").Select() | Out-Null
